$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A41").Value = 46031
$ws.Range("B41").Value = "Faire la documentation"
$ws.Range("D41").Value = 5

$ws.Range("B42").Value = "traduction et implémentation du guide utilisateur"
$ws.Range("D42").Value = 2

$ws.Range("B47").Value = "Journée très longue car faire de la documentation, quelle qu'elle soit pendant 5h c'est long ennuyant et embêtant. Heureusement que Andrei avait besoin de mon aide pour la traduction et l'implémentation du guide utilisateur pour le clavier. Je fûs distrait par l'extérieur plusieurs fois durant la journée pour manque de motivation à faire la documentation."
